$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1101.7059
$ws.Range("I15").Value = 1101.7059
$ws.Range("K15").Value = 3305.1177
$ws.Range("M15").Value = -3136.1177
$ws.Range("H18").Value = 1006.6667
$ws.Range("I18").Value = 1006.6667
$ws.Range("K18").Value = 1006.6667
$ws.Range("M18").Value = -722.6667
$ws.Range("H20").Value = 4432.5
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20460
$ws.Range("H26").Value = 2853
$ws.Range("J26").Value = 2699.5
$ws.Range("L26").Value = 2699.5
$ws.Range("N26").Value = -3387.5
$ws.Range("H35").Value = 4432.5
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20758
$ws.Range("H46").Value = 1333.3334
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -2881
$ws.Range("N46").Value = -6238
$ws.Range("H47").Value = 16250
$ws.Range("I47").Value = 16250
$ws.Range("K47").Value = 16250
$ws.Range("M47").Value = -15278
$ws.Range("H60").Value = 1333.3334
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 2000
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 6000
$ws.Range("M60").Value = -2516
$ws.Range("N60").Value = -6968
$ws.Range("H70").Value = 6199.6
$ws.Range("I70").Value = 5499.5
$ws.Range("J70").Value = 6374.625
$ws.Range("K70").Value = 16498.5
$ws.Range("L70").Value = 19123.875
$ws.Range("M70").Value = -16228.5
$ws.Range("N70").Value = -19663.875
$ws.Range("H73").Value = 6199.6
$ws.Range("I73").Value = 5499.5
$ws.Range("J73").Value = 6374.625
$ws.Range("K73").Value = 16498.5
$ws.Range("L73").Value = 19123.875
$ws.Range("M73").Value = -15562.5
$ws.Range("N73").Value = -20995.875

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 906.1053000000001
$ws.Range("I2").Value = 821.36365
$ws.Range("J2").Value = 1022.625
$ws.Range("K2").Value = 821.36365
$ws.Range("L2").Value = 1022.625
$ws.Range("M2").Value = -708.36365
$ws.Range("N2").Value = -1248.625
$ws.Range("H45").Value = 3112.375
$ws.Range("I45").Value = 2293.1
$ws.Range("K45").Value = 2293.1
$ws.Range("M45").Value = -1916.1
$ws.Range("H63").Value = 10488
$ws.Range("I63").Value = 598.75
$ws.Range("K63").Value = 598.75
$ws.Range("M63").Value = 87.25
$ws.Range("H66").Value = 10488
$ws.Range("I66").Value = 598.75
$ws.Range("K66").Value = 2993.75
$ws.Range("M66").Value = 438.25
$ws.Range("H110").Value = 297.625
$ws.Range("I110").Value = 297.625
$ws.Range("K110").Value = 297.625
$ws.Range("M110").Value = 1747.375
$ws.Range("H116").Value = 906.1053000000001
$ws.Range("I116").Value = 821.36365
$ws.Range("J116").Value = 1022.625
$ws.Range("K116").Value = 821.36365
$ws.Range("L116").Value = 1022.625
$ws.Range("M116").Value = 1472.63635
$ws.Range("N116").Value = -5610.625

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 906.1053000000001
$ws.Range("I3").Value = 821.36365
$ws.Range("J3").Value = 1022.625
$ws.Range("K3").Value = 821.36365
$ws.Range("L3").Value = 1022.625
$ws.Range("M3").Value = -707.36365
$ws.Range("N3").Value = -1250.625
$ws.Range("H94").Value = 1234.5883
$ws.Range("I94").Value = 1289.3572
$ws.Range("J94").Value = 979
$ws.Range("K94").Value = 1289.3572
$ws.Range("L94").Value = 979
$ws.Range("M94").Value = -838.3571999999999
$ws.Range("N94").Value = -1881
$ws.Range("H99").Value = 1732.1538
$ws.Range("I99").Value = 1742.1
$ws.Range("J99").Value = 1699
$ws.Range("K99").Value = 1742.1
$ws.Range("L99").Value = 1699
$ws.Range("M99").Value = -244.0999999999999
$ws.Range("N99").Value = -4695
$ws.Range("H134").Value = 3165.2856
$ws.Range("I134").Value = 3165.2856
$ws.Range("K134").Value = 9495.856800000001
$ws.Range("M134").Value = -6960.856800000001

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5424.579
$ws.Range("I7").Value = 8502.25
$ws.Range("J7").Value = 148.57143
$ws.Range("K7").Value = 8502.25
$ws.Range("L7").Value = 148.57143
$ws.Range("M7").Value = -8389.25
$ws.Range("N7").Value = -374.57143

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 281.2
$ws.Range("J40").Value = 350.125
$ws.Range("L40").Value = 1400.5
$ws.Range("N40").Value = -1538.5
$ws.Range("H61").Value = 138
$ws.Range("I61").Value = 144.14285
$ws.Range("K61").Value = 432.42855
$ws.Range("M61").Value = -217.42855
$ws.Range("H86").Value = 373.6154
$ws.Range("I86").Value = 424.83334
$ws.Range("J86").Value = 329.7143
$ws.Range("K86").Value = 1274.50002
$ws.Range("L86").Value = 989.1428999999999
$ws.Range("M86").Value = -88.50001999999995
$ws.Range("N86").Value = -3361.1429
$ws.Range("H89").Value = 373.6154
$ws.Range("I89").Value = 424.83334
$ws.Range("J89").Value = 329.7143
$ws.Range("K89").Value = 3823.50006
$ws.Range("L89").Value = 2967.4287
$ws.Range("M89").Value = 2104.49994
$ws.Range("N89").Value = -14823.4287
$ws.Range("H117").Value = 2643.4443
$ws.Range("I117").Value = 2000
$ws.Range("J117").Value = 2827.2856
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 8481.856800000001
$ws.Range("M117").Value = -2558
$ws.Range("N117").Value = -15365.8568

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5194.4443
$ws.Range("I113").Value = 2791.6667
$ws.Range("K113").Value = 2791.6667
$ws.Range("M113").Value = -621.6667000000002
$ws.Range("H132").Value = 4012
$ws.Range("I132").Value = 4012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12036
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9506
$ws.Range("N132").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1900.25
$ws.Range("I17").Value = 1101
$ws.Range("J17").Value = 2166.6667
$ws.Range("K17").Value = 1101
$ws.Range("L17").Value = 2166.6667
$ws.Range("M17").Value = -931
$ws.Range("N17").Value = -2506.6667
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H61").Value = 4740.3
$ws.Range("I61").Value = 2769.6
$ws.Range("J61").Value = 6711
$ws.Range("K61").Value = 2769.6
$ws.Range("L61").Value = 6711
$ws.Range("M61").Value = -2567.6
$ws.Range("N61").Value = -7115
$ws.Range("H94").Value = 21000
$ws.Range("J94").Value = 21000
$ws.Range("L94").Value = 21000
$ws.Range("N94").Value = -22352
$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 10000
$ws.Range("N100").Value = -11082
$ws.Range("H113").Value = 4740.3
$ws.Range("I113").Value = 2769.6
$ws.Range("J113").Value = 6711
$ws.Range("K113").Value = 2769.6
$ws.Range("L113").Value = 6711
$ws.Range("M113").Value = -599.5999999999999
$ws.Range("N113").Value = -11051
$ws.Range("H136").Value = 4848
$ws.Range("I136").Value = 3313.5
$ws.Range("K136").Value = 9940.5
$ws.Range("M136").Value = -7390.5
$ws.Range("N58").ClearContents()
$ws.Range("M100").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21622
$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -68112
$ws.Range("H107").Value = 2089.5715
$ws.Range("I107").Value = 2321.1667
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 6963.500100000001
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = -5043.500100000001
$ws.Range("N107").Value = -5940
$ws.Range("H113").Value = 843.625
$ws.Range("I113").Value = 687.25
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2061.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 108.25
$ws.Range("N113").Value = -7340
